$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells matching the existing header style (copy format from AC1)
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill team record (Wins/Losses/Ties) for every data row
$lastRow = 55
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 89
    $ws.Cells.Item($r, 31).Value = 73
    $ws.Cells.Item($r, 32).Value = 0
}
